$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '66.884.54'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.90%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.514.73'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '607.49'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '147.92'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.05%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.513.59'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.60%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E9').Value = '  -1.66%  '
$ws.Range('E10').Value = '  -1.09%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '8.01'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +5.53%  '
$ws.Range('E12').Value = '  -1.93%  '
$ws.Range('E13').Value = '  +0.90%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.109.33'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.70%  '
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.509.73'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '66.971.79'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.91%  '
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.77'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +8.68%  '
$ws.Range('E20').Value = '  -0.71%  '
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '438.30'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.94%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.610'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.74%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '79.47'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.21%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.653.82'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  -3.89%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.79'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.70%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.32'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.22%  '
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('E31').Value = '  -1.50%  '
$ws.Range('E32').Value = '  -1.68%  '
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '25.56'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.98'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.56%  '
$ws.Range('E36').Value = '  -2.18%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '8.05'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.52%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '174.11'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0897'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('E42').Value = '  -0.31%  '
$ws.Range('E43').Value = '  -10.32%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.894'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '46.17'
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '27.92'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -7.81%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.27'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.47'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.93%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.46'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.02%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.993'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.40%  '
$ws.Range('E51').Value = '  -0.98%  '
